# Quarterly indexing esoteric bug-fix operation
#
# Column A (rows 2-63) stores date serial numbers that should represent the
# 15th of the month following the previously-recorded month (instead of the
# 1st of that month). This shifts every date by one month, landing on the
# 15th, which nets out to +44 or +45 days depending on the month lengths.
#
# New serial values for rows 2..63 (row N corresponds to $newValues[N-2]):
$newValues = @(
    25614, 25614, 25614, 25614, 25614, 25614, 25614, 25614, 25614, 25614,
    40313, 40405, 40497, 40589, 40678, 40770, 40862, 40954, 41044, 41136,
    41228, 41320, 41409, 41501, 41593, 41685, 41774, 41866, 41958, 42050,
    42139, 42231, 42323, 42415, 42505, 42597, 42689, 42781, 42870, 42962,
    43054, 43146, 43235, 43327, 43419, 43511, 43600, 43692, 43784, 43876,
    43966, 44058, 44150, 44242, 44331, 44423, 44515, 44607, 44696, 44788,
    44880, 44972
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value2 = $newValues[$i]
}
